$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$after = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $after)
$ws.Name = "Sheet4"

# Column A: sequence 0..23
$ws.Range("A1").Value = 0
for ($r = 2; $r -le 24; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=A$prev+1"
}

# Column B: A*8
for ($r = 1; $r -le 24; $r++) {
    $ws.Range("B$r").Formula = "=A$r*8"
}

# Column C: DEC2HEX(B, 2), right aligned (matches existing style index 7)
for ($r = 1; $r -le 24; $r++) {
    $ws.Range("C$r").Formula = "=DEC2HEX(B$r, 2)"
}
$ws.Range("C1:C24").HorizontalAlignment = -4152

# Column E: page numbers (first 8 numeric-looking, last 3 textual -> shared strings)
$ws.Range("E1").Value = "2108"
$ws.Range("E2").Value = "6908"
$ws.Range("E3").Value = "7108"
$ws.Range("E4").Value = "7908"
$ws.Range("E5").Value = "8108"
$ws.Range("E6").Value = "8908"
$ws.Range("E7").Value = "9108"
$ws.Range("E8").Value = "9908"
$ws.Range("E9").Value = "A108"
$ws.Range("E10").Value = "A908"
$ws.Range("E11").Value = "B108"
$ws.Range("E1:E11").HorizontalAlignment = -4131

# Column widths (approx default width, matches target 9.140625)
$ws.Columns.Item(3).ColumnWidth = 9
$ws.Columns.Item(5).ColumnWidth = 9

$ws.Calculate()

# Selection / view state
[void]$ws.Range("G16").Select()

Write-Host "Sheet4 added with" $wb.Worksheets.Count "sheets total"
